$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 16 mirrors the pattern of existing data rows (A=index, B=label, C:M=values).
# Copy row 15's formatting for column A down to A16 so it keeps the same bold/
# bordered/centered style (style index 1) without minting a new cell style.
$ws.Range("A15").Copy($ws.Range("A16"))

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"

$ws.Range("C16").Value = 0.9917500818181104
$ws.Range("D16").Value = 1.00398004565772
$ws.Range("E16").Value = 0.9983221536969686
$ws.Range("F16").Value = 0.9917500818181104
$ws.Range("G16").Value = 0.9993463988198579
$ws.Range("H16").Value = 0.9929919182923702
$ws.Range("I16").Value = 0.9952941176470588
$ws.Range("J16").Value = 1.00398004565772
$ws.Range("K16").Value = 1.001151099677344
$ws.Range("L16").Value = 0.9964505907477272
$ws.Range("M16").Value = 0.9969474526553476
